# Fix import dan menghilangkan total potongan dan total tunjangan
#
# The sheet had two "total" columns that were really just duplicates of the
# per-item 45.000 figure (they weren't genuine SUM() totals, just stray
# 675000 / 135000 constants left over from a copy/paste): column V
# ("TOTAL POTONGAN") and the last column, originally Z ("TOTAL TUNJANGAN").
# Both get removed entirely (entire-column delete), which shifts every
# column to their right one step to the left and drops the now-unused
# "TOTAL POTONGAN" / "TOTAL TUNJANGAN" shared strings automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "TOTAL POTONGAN" column (was column V).
$ws.Range("V1").EntireColumn.Delete() | Out-Null

# After the delete above, everything shifted left by one, so the old
# "TOTAL TUNJANGAN" column (was Z) is now Y. Delete it too.
$ws.Range("Y1").EntireColumn.Delete() | Out-Null

# Leave the selection where the author's session ended up.
$ws.Range("X2").Select() | Out-Null
